# Update the player roster table (columns A:C, rows 2-19) on the active sheet
# to reflect the refreshed player/position/team data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Damian Lillard",      "PG",       "Milwaukee Bucks"),
    @("Cade Cunningham",     "PG,SG",    "Detroit Pistons"),
    @("LaMelo Ball",         "PG,SG",    "Charlotte Hornets"),
    @("Devin Vassell",       "SG,SF",    "San Antonio Spurs"),
    @("Ausar Thompson",      "SF,PF",    "Detroit Pistons"),
    @("Isaiah Hartenstein",  "C",        "Oklahoma City Thunder"),
    @("Onyeka Okongwu",      "PF,C",     "Atlanta Hawks"),
    @("Cam Whitmore",        "SF,PF",    "Houston Rockets"),
    @("Naz Reid",            "PF,C",     "Minnesota Timberwolves"),
    @("T.J. McConnell",      "PG",       "Indiana Pacers"),
    @("Derrick White",       "PG,SG",    "Boston Celtics"),
    @("Mike Conley",         "PG",       "Minnesota Timberwolves"),
    @("Deandre Ayton",       "C",        "Portland Trail Blazers"),
    @("Malik Monk",          "PG,SG,SF", "Sacramento Kings"),
    @("Coby White",          "PG,SG",    "Chicago Bulls"),
    @("Anthony Davis",       "PF,C",     "Dallas Mavericks"),
    @("Collin Sexton",       "PG,SG",    "Utah Jazz"),
    @("Julius Randle",       "PF,C",     "Minnesota Timberwolves")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
